# Generate Report for Archive
#
# 1) Update the "Status" value everywhere it reads "Ready for handoff" to
#    "In Translation" (Overview sheet columns E/F, and the Status column
#    ("C") on the zh-cn / de-de sheets).
# 2) Narrow the Status-related columns (Overview!E:F, zh-cn!C, de-de!C)
#    from their current width down to ~13.41 characters.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
for ($si = 1; $si -le $sheetCount; $si++) {
    $ws = $wb.Worksheets.Item($si)
    $used = $ws.UsedRange
    $nrows = $used.Rows.Count
    $ncols = $used.Columns.Count
    for ($r = 1; $r -le $nrows; $r++) {
        for ($c = 1; $c -le $ncols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            # NOTE: compare with the literal on the left -- some cell values
            # (e.g. boolean-looking "True"/"False" text) come back as a
            # native Boolean from .Text/.Value, and "-eq" coerces its right
            # operand to the type of the *left* operand. Putting the string
            # literal first keeps this a plain string comparison.
            if ("Ready for handoff" -eq [string]$cell.Text) {
                $cell.Value = "In Translation"
            }
        }
    }
}

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 12.5
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 12.5
